$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [string][char]0x2083

$ws.Range("D2").Value = "25.864.24"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.631.15"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.506"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "1.662.75"
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("D13").Value = "1.856.18"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.545"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0" + $sub3 + "0755"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "25.875.50"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  +2.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("D37").Value = "1.136.97"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.796"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "1.765.50"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  +3.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0527"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.94%  "
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.93%  "
